$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day-1 (rows 2-97): explicit overrides for Notified/Actual Production (MW) found in the
# commit's updated forecast (GESS added to the portfolio). Rows not listed keep their
# existing Notified(B)/Actual(C) values.
$day1B = @{2=0; 3=0; 4=0; 5=0; 6=0; 7=0; 8=0; 9=0; 10=0; 11=0; 12=0; 13=0; 22=2; 23=2; 24=2; 25=3; 26=18; 27=23; 28=36; 29=54; 30=205; 31=255; 32=313; 33=375; 34=625; 35=709; 36=787; 37=865; 38=1038; 39=1110; 40=1166; 41=1221; 42=1344; 43=1385; 44=1417; 45=1449; 46=1469; 47=1480; 48=1482; 49=1490; 50=1457; 51=1445; 52=1433; 53=1410; 54=1351; 55=1314; 56=1276; 57=1240; 58=1132; 59=1084; 60=1034; 61=973; 62=841; 63=803; 64=699; 65=621; 66=453; 67=371; 68=310; 69=256; 70=119; 71=84; 72=52; 73=39; 74=14; 75=12; 76=9; 77=8; 78=9; 79=8; 80=7; 81=7; 82=1; 83=1; 84=0; 85=0; 94=0; 95=0; 96=0; 97=0}
$day1C = @{27=1; 28=13; 29=39; 30=79; 31=130; 32=194; 33=276; 34=360; 35=443; 36=513; 37=600; 38=0; 39=0; 40=0; 41=0; 42=0; 43=0; 44=0; 45=0; 46=0; 47=0; 48=0; 49=0; 50=0; 51=0; 52=0; 53=0; 54=0; 55=0; 56=0; 57=0; 58=0; 59=0; 60=0; 61=0; 62=0; 63=0; 64=0; 65=0; 66=0; 67=0; 68=0; 69=0; 70=0; 71=0; 72=0; 73=0; 74=0; 75=0}

for ($r = 2; $r -le 193; $r++) {
    # Shift the Timestamp (column A) forward by 9 days: 22/23 Sep 2025 -> 1/2 Oct 2025.
    $oldA = $ws.Cells.Item($r, 1).Value2()
    $ws.Cells.Item($r, 1).Value = $oldA + 9

    if ($r -le 97) {
        $q = $r - 1
        $lookup = "01.10.2025" + $q
        if ($day1B.ContainsKey($r)) {
            $ws.Cells.Item($r, 2).Value = $day1B[$r]
        }
        if ($day1C.ContainsKey($r)) {
            $ws.Cells.Item($r, 3).Value = $day1C[$r]
        }
    } else {
        $q = $r - 97
        $lookup = "02.10.2025" + $q
        # Day-2 (rows 98-193): no production yet forecast/observed -> zero out both columns.
        $ws.Cells.Item($r, 2).Value = 0
        $ws.Cells.Item($r, 3).Value = 0
    }

    $ws.Cells.Item($r, 5).Value = $lookup
}
